$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "2025-04-29 04:45:47"
$ws.Range("B44").Value = 146
